{"js": "// Find the target run's text and split it into two runs:\n//  - \"Chapter 4: The evolution of the retinol metabolism \" (unchanged formatting)\n//  - \"and its role in the origin of vision\" (same formatting + yellow highlight)\nconst fullText =\n  \"Chapter 4: The evolution of the retinol metabolism and its role in the origin of vision\";\nconst highlightText = \"and its role in the origin of vision\";\n\nconst results = context.document.body.search(fullText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph text not found\");\n}\n\nconst target = results.items[0];\n\n// Locate the sub-range that needs the yellow highlight, searched within the\n// matched range so we keep the split scoped to this paragraph only.\nconst subResults = target.search(highlightText, { matchCase: true });\nsubResults.load(\"items\");\nawait context.sync();\n\nif (subResults.items.length === 0) {\n  throw new Error(\"Highlight sub-text not found\");\n}\n\nconst highlightRange = subResults.items[0];\nhighlightRange.font.highlightColor = \"Yellow\";\n\nawait context.sync();\n", "ps1": "# Split the Chapter 4 title run into two runs so the second half\n# (\"and its role in the origin of vision\") can carry a yellow highlight,\n# matching the other chapter titles' highlighted continuations.\n$d = $word.ActiveDocument\n\n$fullText = \"Chapter 4: The evolution of the retinol metabolism and its role in the origin of vision\"\n$highlightText = \"and its role in the origin of vision\"\n\n# Locate the unique target paragraph's text first, scoping everything that\n# follows to this exact Range so we don't touch the similarly-worded\n# \"Chapter 5\" paragraph earlier in the document.\n$target = $d.Content\n$target.Find.ClearFormatting()\n$target.Find.Text = $fullText\n$target.Find.MatchCase = $true\n$target.Find.MatchWholeWord = $false\n$target.Find.Wrap = 0\nif (-not $target.Find.Execute()) {\n    throw \"Target chapter title text not found\"\n}\n\n# Now search only within that matched range for the trailing phrase that\n# needs the new highlighted run.\n$highlightRange = $target.Duplicate\n$highlightRange.Find.ClearFormatting()\n$highlightRange.Find.Text = $highlightText\n$highlightRange.Find.MatchCase = $true\n$highlightRange.Find.Wrap = 0\nif (-not $highlightRange.Find.Execute()) {\n    throw \"Highlight sub-text not found\"\n}\n\n# Setting Font.HighlightColorIndex on this sub-range splits the original\n# single run into two runs, the second one carrying the highlight -\n# exactly mirroring the target OOXML diff.\n$highlightRange.Font.HighlightColorIndex = 7  # wdYellow\n\n$d.Save()\n"}
